$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Cells.Item(4, 8).Value = 224.85715
$ws.Cells.Item(4, 9).Value = 224.85715
$ws.Cells.Item(4, 11).Value = 224.85715
$ws.Cells.Item(4, 13).Value = -110.85715

# Row 64
$ws.Cells.Item(64, 8).Value = 2400
$ws.Cells.Item(64, 9).Value = 2400
$ws.Cells.Item(64, 11).Value = 2400
$ws.Cells.Item(64, 13).Value = -2152

# Row 67
$ws.Cells.Item(67, 8).Value = 2400
$ws.Cells.Item(67, 9).Value = 2400
$ws.Cells.Item(67, 11).Value = 2400
$ws.Cells.Item(67, 13).Value = -1542

# Row 86
$ws.Cells.Item(86, 8).Value = 75146.336
$ws.Cells.Item(86, 9).Value = 4082
$ws.Cells.Item(86, 10).Value = 92912.414
$ws.Cells.Item(86, 11).Value = 4082
$ws.Cells.Item(86, 12).Value = 92912.414
$ws.Cells.Item(86, 13).Value = -2959
$ws.Cells.Item(86, 14).Value = -95158.414

# Row 89
$ws.Cells.Item(89, 8).Value = 75146.336
$ws.Cells.Item(89, 9).Value = 4082
$ws.Cells.Item(89, 10).Value = 92912.414
$ws.Cells.Item(89, 11).Value = 20410
$ws.Cells.Item(89, 12).Value = 464562.07
$ws.Cells.Item(89, 13).Value = -14794
$ws.Cells.Item(89, 14).Value = -475794.07

# Row 94
$ws.Cells.Item(94, 8).Value = 4019.7144
$ws.Cells.Item(94, 9).Value = 4019.7144
$ws.Cells.Item(94, 10).Value = 0
$ws.Cells.Item(94, 11).Value = 4019.7144
$ws.Cells.Item(94, 12).Value = 0
$ws.Cells.Item(94, 13).Value = ""
$ws.Cells.Item(94, 14).Value = -3568.7144


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 6
$ws.Cells.Item(6, 8).Value = 59637
$ws.Cells.Item(6, 9).Value = 63086.293
$ws.Cells.Item(6, 11).Value = 63086.293
$ws.Cells.Item(6, 13).Value = -62913.293

# Row 32
$ws.Cells.Item(32, 8).Value = 3422.3684
$ws.Cells.Item(32, 9).Value = 2135
$ws.Cells.Item(32, 11).Value = 2135
$ws.Cells.Item(32, 13).Value = -1848

# Row 74
$ws.Cells.Item(74, 8).Value = 1239.125
$ws.Cells.Item(74, 9).Value = 1130.4286
$ws.Cells.Item(74, 10).Value = 2000
$ws.Cells.Item(74, 11).Value = 1130.4286
$ws.Cells.Item(74, 12).Value = 2000
$ws.Cells.Item(74, 13).Value = -256.4286
$ws.Cells.Item(74, 14).Value = -3748

# Row 77
$ws.Cells.Item(77, 8).Value = 1239.125
$ws.Cells.Item(77, 9).Value = 1130.4286
$ws.Cells.Item(77, 10).Value = 2000
$ws.Cells.Item(77, 11).Value = 5652.143
$ws.Cells.Item(77, 12).Value = 10000
$ws.Cells.Item(77, 13).Value = -1284.143
$ws.Cells.Item(77, 14).Value = -18736

# Row 92
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 12).Value = ""
$ws.Cells.Item(92, 14).Value = 0

# Row 95
$ws.Cells.Item(95, 8).Value = 50000
$ws.Cells.Item(95, 10).Value = 50000
$ws.Cells.Item(95, 12).Value = 50000
$ws.Cells.Item(95, 14).Value = -55492

# Row 97
$ws.Cells.Item(97, 8).Value = 436.8889
$ws.Cells.Item(97, 9).Value = 379
$ws.Cells.Item(97, 11).Value = 379
$ws.Cells.Item(97, 13).Value = 117


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Cells.Item(5, 8).Value = 122.52941
$ws.Cells.Item(5, 10).Value = 101.666664
$ws.Cells.Item(5, 12).Value = 101.666664
$ws.Cells.Item(5, 14).Value = -327.666664

# Row 7
$ws.Cells.Item(7, 8).Value = 590.3333
$ws.Cells.Item(7, 9).Value = 291.8
$ws.Cells.Item(7, 10).Value = 963.5
$ws.Cells.Item(7, 11).Value = 291.8
$ws.Cells.Item(7, 12).Value = 963.5
$ws.Cells.Item(7, 13).Value = -178.8
$ws.Cells.Item(7, 14).Value = -1189.5

# Row 86
$ws.Cells.Item(86, 8).Value = 3312
$ws.Cells.Item(86, 9).Value = 3499.8
$ws.Cells.Item(86, 10).Value = 2999
$ws.Cells.Item(86, 11).Value = 3499.8
$ws.Cells.Item(86, 12).Value = 2999
$ws.Cells.Item(86, 13).Value = -2376.8
$ws.Cells.Item(86, 14).Value = -5245

# Row 89
$ws.Cells.Item(89, 8).Value = 3312
$ws.Cells.Item(89, 9).Value = 3499.8
$ws.Cells.Item(89, 10).Value = 2999
$ws.Cells.Item(89, 11).Value = 17499
$ws.Cells.Item(89, 12).Value = 14995
$ws.Cells.Item(89, 13).Value = -11883
$ws.Cells.Item(89, 14).Value = -26227

# Row 105
$ws.Cells.Item(105, 8).Value = 5636.091
$ws.Cells.Item(105, 9).Value = 2888.0625
$ws.Cells.Item(105, 10).Value = 12964.167
$ws.Cells.Item(105, 11).Value = 2888.0625
$ws.Cells.Item(105, 12).Value = 12964.167
$ws.Cells.Item(105, 13).Value = -1141.0625
$ws.Cells.Item(105, 14).Value = -16458.167


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 3
$ws.Cells.Item(3, 8).Value = 219.83333
$ws.Cells.Item(3, 9).Value = 219.83333
$ws.Cells.Item(3, 11).Value = 219.83333
$ws.Cells.Item(3, 13).Value = -106.83333

# Row 16
$ws.Cells.Item(16, 8).Value = 1944
$ws.Cells.Item(16, 9).Value = 1944
$ws.Cells.Item(16, 11).Value = 1944
$ws.Cells.Item(16, 13).Value = -1657

# Row 32
$ws.Cells.Item(32, 8).Value = 790.375
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 14).Value = ""

# Row 62
$ws.Cells.Item(62, 8).Value = 32699.867
$ws.Cells.Item(62, 9).Value = 5750.1113
$ws.Cells.Item(62, 10).Value = 73124.5
$ws.Cells.Item(62, 11).Value = 5750.1113
$ws.Cells.Item(62, 12).Value = 73124.5
$ws.Cells.Item(62, 13).Value = -5126.1113
$ws.Cells.Item(62, 14).Value = -74372.5

# Row 65
$ws.Cells.Item(65, 8).Value = 32699.867
$ws.Cells.Item(65, 9).Value = 5750.1113
$ws.Cells.Item(65, 10).Value = 73124.5
$ws.Cells.Item(65, 11).Value = 28750.5565
$ws.Cells.Item(65, 12).Value = 365622.5
$ws.Cells.Item(65, 13).Value = -25630.5565
$ws.Cells.Item(65, 14).Value = -371862.5

# Row 93
$ws.Cells.Item(93, 8).Value = 1500
$ws.Cells.Item(93, 9).Value = 1500
$ws.Cells.Item(93, 11).Value = 1500
$ws.Cells.Item(93, 13).Value = 372

# Row 97
$ws.Cells.Item(97, 8).Value = 60000
$ws.Cells.Item(97, 10).Value = 60000
$ws.Cells.Item(97, 12).Value = 60000
$ws.Cells.Item(97, 14).Value = -61982

# Row 113
$ws.Cells.Item(113, 8).Value = 1944
$ws.Cells.Item(113, 9).Value = 1944
$ws.Cells.Item(113, 11).Value = 1944
$ws.Cells.Item(113, 13).Value = 226


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Cells.Item(14, 8).Value = 1807.25
$ws.Cells.Item(14, 9).Value = 1807.25
$ws.Cells.Item(14, 11).Value = 5421.75
$ws.Cells.Item(14, 13).Value = -5248.75

# Row 49
$ws.Cells.Item(49, 8).Value = 232.22223
$ws.Cells.Item(49, 9).Value = 250
$ws.Cells.Item(49, 10).Value = 196.66667
$ws.Cells.Item(49, 11).Value = 750
$ws.Cells.Item(49, 12).Value = 590.00001
$ws.Cells.Item(49, 13).Value = -594
$ws.Cells.Item(49, 14).Value = -902.00001

# Row 138
$ws.Cells.Item(138, 8).Value = 4000
$ws.Cells.Item(138, 9).Value = 5000
$ws.Cells.Item(138, 11).Value = 15000
$ws.Cells.Item(138, 13).Value = -9860


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Cells.Item(126, 8).Value = 4141.3335
$ws.Cells.Item(126, 9).Value = 3068
$ws.Cells.Item(126, 11).Value = 9204
$ws.Cells.Item(126, 13).Value = -6734


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 6
$ws.Cells.Item(6, 8).Value = 43000
$ws.Cells.Item(6, 10).Value = 43000
$ws.Cells.Item(6, 12).Value = 43000
$ws.Cells.Item(6, 14).Value = -43224

# Row 16
$ws.Cells.Item(16, 8).Value = 39974
$ws.Cells.Item(16, 9).Value = 43298.668
$ws.Cells.Item(16, 11).Value = 43298.668
$ws.Cells.Item(16, 13).Value = -43128.668

# Row 93
$ws.Cells.Item(93, 8).Value = 8596.4
$ws.Cells.Item(93, 9).Value = 10997.333
$ws.Cells.Item(93, 11).Value = 10997.333
$ws.Cells.Item(93, 13).Value = -9749.333000000001

# Row 99
$ws.Cells.Item(99, 8).Value = 55000
$ws.Cells.Item(99, 10).Value = 55000
$ws.Cells.Item(99, 12).Value = 55000
$ws.Cells.Item(99, 14).Value = -60990

# Row 101
$ws.Cells.Item(101, 8).Value = 20199.5
$ws.Cells.Item(101, 10).Value = 20199.5
$ws.Cells.Item(101, 12).Value = 20199.5
$ws.Cells.Item(101, 14).Value = -26689.5


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Cells.Item(2, 8).Value = 5904.5713
$ws.Cells.Item(2, 9).Value = 3767.3333
$ws.Cells.Item(2, 11).Value = 3767.3333
$ws.Cells.Item(2, 13).Value = -3655.3333

# Row 69
$ws.Cells.Item(69, 8).Value = 5998.8335
$ws.Cells.Item(69, 10).Value = 5998.8335
$ws.Cells.Item(69, 12).Value = 5998.8335
$ws.Cells.Item(69, 14).Value = -7496.8335

# Row 72
$ws.Cells.Item(72, 8).Value = 5998.8335
$ws.Cells.Item(72, 10).Value = 5998.8335
$ws.Cells.Item(72, 12).Value = 17996.5005
$ws.Cells.Item(72, 14).Value = -25484.5005

# Row 81
$ws.Cells.Item(81, 8).Value = 1990.3636
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 12).Value = 0
$ws.Cells.Item(81, 14).Value = ""

# Row 84
$ws.Cells.Item(84, 8).Value = 1990.3636
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 14).Value = ""

# Row 133
$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 10).Value = 0
$ws.Cells.Item(133, 12).Value = ""
$ws.Cells.Item(133, 14).Value = 0

